$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row with mixed (rich-text) formatting: "Bold and Italic mixed"
$cell = $ws.Range("A5")
$cell.Value = "Bold and Italic mixed"

# "Bold" -> characters 1-4
$cell.Characters(1, 4).Font.Bold = $true
# "Italic" -> characters 10-6 ("Bold and Italic mixed"; "Italic" starts at position 10)
$cell.Characters(10, 6).Font.Italic = $true

$ws.Range("A5").Select()
